# Applies the attendance_reports sync edit described in the commit:
# "Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-06 09:16:25"
#
# Changes:
#  1. Column I (9) width 10 -> 14
#  2. "Recorded By" (column G) name ordering reshuffled on many rows
#  3. Class Statistics swap: Missing Sessions (L7) <-> Pending Sessions (L8) counts
#  4. P18:Q18, P19:Q19, P20:Q20 value swaps
#  5. Status "Pending" -> "Not Recorded" on rows 105, 131, 157

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column I width change (10 -> 14). Excel's ColumnWidth setter rounds to the
#    workbook's default-font pixel grid, so 13.15 is the input that resolves to
#    a stored width of exactly 14 characters.
$ws.Columns.Item(9).ColumnWidth = 13.15

# 2. "Recorded By" (column G) text reorderings.
#    "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System"
$sysSysBackupRows = @(2, 28, 54)
foreach ($r in $sysSysBackupRows) {
    $ws.Range("G$r").Value = "backup@backdoor.com, system, System"
}

#    "System, backup@backdoor.com" -> "backup@backdoor.com, System"
$sysBackupRows = @(4, 5, 8, 30, 31, 34, 56, 57, 60, 80, 81, 82, 106, 107, 108, 132, 133, 134)
foreach ($r in $sysBackupRows) {
    $ws.Range("G$r").Value = "backup@backdoor.com, System"
}

#    "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$dnasrSystemRows = @(11, 17, 37, 43, 63, 69, 93, 94, 96, 119, 120, 122, 145, 146, 148)
foreach ($r in $dnasrSystemRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# 3. Class Statistics: Missing Sessions / Pending Sessions counts swap.
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 0

# 4. Per-group P/Q swaps for rows 18, 19, 20.
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 0
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 0
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 0

# 5. Status text updates.
$pendingToNotRecordedRows = @(105, 131, 157)
foreach ($r in $pendingToNotRecordedRows) {
    $ws.Range("I$r").Value = "Not Recorded"
}

$wb.Save()
